$d = $word.ActiveDocument

$d.Content.Find.Execute("82×18=", $true, $false, $false, $false, $false, $true, 1, $false, "43×60=", 2) | Out-Null
$d.Content.Find.Execute("99×30=", $true, $false, $false, $false, $false, $true, 1, $false, "88×42=", 2) | Out-Null
$d.Content.Find.Execute("25×23=", $true, $false, $false, $false, $false, $true, 1, $false, "41×87=", 2) | Out-Null
$d.Content.Find.Execute("64×38=", $true, $false, $false, $false, $false, $true, 1, $false, "25×65=", 2) | Out-Null
$d.Content.Find.Execute("22×73=", $true, $false, $false, $false, $false, $true, 1, $false, "28×27=", 2) | Out-Null
$d.Content.Find.Execute("71×92=", $true, $false, $false, $false, $false, $true, 1, $false, "38×93=", 2) | Out-Null
$d.Content.Find.Execute("75×57=", $true, $false, $false, $false, $false, $true, 1, $false, "31×53=", 2) | Out-Null
$d.Content.Find.Execute("38×52=", $true, $false, $false, $false, $false, $true, 1, $false, "39×19=", 2) | Out-Null
$d.Content.Find.Execute("56×52=", $true, $false, $false, $false, $false, $true, 1, $false, "15×93=", 2) | Out-Null
$d.Content.Find.Execute("55×99=", $true, $false, $false, $false, $false, $true, 1, $false, "51×61=", 2) | Out-Null
$d.Content.Find.Execute("46×89=", $true, $false, $false, $false, $false, $true, 1, $false, "26×93=", 2) | Out-Null
$d.Content.Find.Execute("79×55=", $true, $false, $false, $false, $false, $true, 1, $false, "93×36=", 2) | Out-Null
$d.Content.Find.Execute("92×21=", $true, $false, $false, $false, $false, $true, 1, $false, "13×74=", 2) | Out-Null
$d.Content.Find.Execute("29×35=", $true, $false, $false, $false, $false, $true, 1, $false, "65×30=", 2) | Out-Null
$d.Content.Find.Execute("36×56=", $true, $false, $false, $false, $false, $true, 1, $false, "58×85=", 2) | Out-Null
$d.Content.Find.Execute("27×56=", $true, $false, $false, $false, $false, $true, 1, $false, "89×81=", 2) | Out-Null
$d.Content.Find.Execute("37×22=", $true, $false, $false, $false, $false, $true, 1, $false, "63×13=", 2) | Out-Null
$d.Content.Find.Execute("41×27=", $true, $false, $false, $false, $false, $true, 1, $false, "83×52=", 2) | Out-Null
$d.Content.Find.Execute("99×99=", $true, $false, $false, $false, $false, $true, 1, $false, "47×74=", 2) | Out-Null
$d.Content.Find.Execute("31×34=", $true, $false, $false, $false, $false, $true, 1, $false, "15×50=", 2) | Out-Null
$d.Content.Find.Execute("15×62=", $true, $false, $false, $false, $false, $true, 1, $false, "20×81=", 2) | Out-Null
$d.Content.Find.Execute("56×50=", $true, $false, $false, $false, $false, $true, 1, $false, "92×73=", 2) | Out-Null
$d.Content.Find.Execute("25×94=", $true, $false, $false, $false, $false, $true, 1, $false, "59×14=", 2) | Out-Null
$d.Content.Find.Execute("29×79=", $true, $false, $false, $false, $false, $true, 1, $false, "30×47=", 2) | Out-Null
$d.Content.Find.Execute("96×65=", $true, $false, $false, $false, $false, $true, 1, $false, "42×35=", 2) | Out-Null
